$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying the existing "2021-Q2" sheet
#    (same column layout/styles as other quarterly sheets), inserted right
#    before "2021-Q2" so the tab order becomes:
#    总计, 2022-Q3, 2021-Q2, 2020-Q4
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q2")
$src.Copy($src)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q3"

# Columns B (fund code) and D:G (percentages/amounts stored as text, like the
# other quarterly sheets) must stay text so leading zeros / exact digits are
# preserved instead of being parsed as numbers.
$ws.Range("B2:B6").NumberFormat = "@"
$ws.Range("D2:G6").NumberFormat = "@"

# This sheet's "D" header reads "基金规模" (fund size) instead of the older
# "基金金额" wording used on the 2021-Q2/2020-Q4 sheets.
$ws.Range("D1").Value = "基金规模"

# Extend column A's number style (from the copied row 2) down through row 6.
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "005698"
$ws.Range("C2").Value = "华夏全球科技先锋混合（QDII）"
$ws.Range("D2").Value = "0.59"
$ws.Range("E2").Value = "86.79"
$ws.Range("F2").Value = "2.39"
$ws.Range("G2").Value = "0.0141"
$ws.Range("H2").Value = 10

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "096001"
$ws.Range("C3").Value = "大成标普500等权重指数（QDII）人民币"
$ws.Range("D3").Value = "3.08"
$ws.Range("E3").Value = "93.16"
$ws.Range("F3").Value = "0.21"
$ws.Range("G3").Value = "0.0065"
$ws.Range("H3").Value = 10

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "013404"
$ws.Range("C4").Value = "大成标普500等权重指数（QDII）美元"
$ws.Range("D4").Value = "3.08"
$ws.Range("E4").Value = "93.16"
$ws.Range("F4").Value = "0.21"
$ws.Range("G4").Value = "0.0065"
$ws.Range("H4").Value = 10

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "011706"
$ws.Range("C5").Value = "长信美国标准普尔100等权重指数增强（QDII）美元"
$ws.Range("D5").Value = "0.39"
$ws.Range("E5").Value = "82.64"
$ws.Range("F5").Value = "0.91"
$ws.Range("G5").Value = "0.0035"
$ws.Range("H5").Value = 3

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "519981"
$ws.Range("C6").Value = "长信美国标准普尔100等权重指数增强（QDII）人民币"
$ws.Range("D6").Value = "0.39"
$ws.Range("E6").Value = "82.64"
$ws.Range("F6").Value = "0.91"
$ws.Range("G6").Value = "0.0035"
$ws.Range("H6").Value = 3

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 above
#    the existing 2021-Q2 row (which pushes 2021-Q2/2020-Q4 down by one row),
#    then fix up the running index in column A.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.03

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# ---------------------------------------------------------------------------
# 3. Restore "2020-Q4" as the active/selected sheet (unchanged from before).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
